$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells whose content moved elsewhere / is removed ---
$ws.Range("D1").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("D2").ClearContents()

# --- Update cells whose text changes in place (string already exists in table) ---
$ws.Range("B1").Value = "Needs tweaking."
$ws.Range("C1").Value = "Shake vigorously with ice. Strain into a {gcocktail glass} and garnish with lime twist."
$ws.Range("C2").Value = "Shake with ice for a long time. Strain into a chilled {gcocktail glass}."
$ws.Range("B3").Value = "What to do about foreign characters like ç? What about common ingredients like sugar?"
$ws.Range("C3").Value = "Muddle lime and sugar in a {glowball glass} until the lime is juiced. Fill to brim with crushed ice and add cachaça. Garnish with sugar cane."
$ws.Range("D3").Value = "Proper Brazilian Caipirinhas are unmeasured and the glass is simply filled with cachaça after the ice is added."

# --- Write new per-ingredient cells (column-major I, E, F, G, H to match authoring order) ---
$ws.Range("I1").Value = "twist of lime"
$ws.Range("I2").Value = "{q1-2} {udashes} orange bitters"

$ws.Range("E1").Value = "gin | {q3/4}{uoz} gin"
$ws.Range("E2").Value = "dry gin | {q1.5}{uoz} dry gin"
$ws.Range("E3").Value = "lime | {q1} {ulime}, cut into eighths"

$ws.Range("F1").Value = "green chartreuse | {q3/4}{uoz} green Chartreuse"
$ws.Range("F2").Value = "dry vermouth | {q1/2}{uoz} dry vermouth"
$ws.Range("F3").Value = "sugar | {q1}{utsp} sugar"

$ws.Range("G1").Value = "maraschino liqueur | {q3/4}{uoz} maraschino liqueur"
$ws.Range("G2").Value = "green chartreuse | {q1/2}{uoz} green Chartreuse"
$ws.Range("G3").Value = "cachaca | {q2}{uoz} cachaça"

$ws.Range("H1").Value = "lime juice | {q3/4}{uoz} lime juice"
$ws.Range("H2").Value = "absinthe | {q1}{utsp} absinthe"

$ws.Range("I3").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1

Write-Host "content done"
